$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 463, shifting existing row 463 (and all following rows) down by one.
$ws.Rows("463:463").Insert()

# Populate the newly inserted row 463 with the new market record.
$ws.Range("A463").Value = 11
$ws.Range("B463").Value = "Vega Monumental Concepción"
$ws.Range("C463").Value = "Bíobío"
$ws.Range("D463").Value = 45218
$ws.Range("E463").Value = 8
$ws.Range("F463").Value = "Fruta"
$ws.Range("G463").Value = 100102
$ws.Range("H463").Value = "Cítricos"
$ws.Range("I463").Value = 100102005
$ws.Range("J463").Value = "Naranja"
$ws.Range("K463").Value = "Navel Late"
$ws.Range("L463").Value = "Primera"
$ws.Range("M463").Value = 300
$ws.Range("N463").Value = 9000
$ws.Range("O463").Value = 9000
$ws.Range("P463").Value = 9000
$ws.Range("Q463").Value = "$/bandeja 15 kilos granel"
$ws.Range("R463").Value = "Región de O'Higgins"
$ws.Range("S463").Value = 600
$ws.Range("T463").Value = 15
